# Applies updated cryptocurrency price/volume data to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value looks like a plain number need to be
# forced to Text format first, otherwise Excel will silently turn them into
# numeric cells instead of preserving the original text-valued display.
$textFormatCells = @("D4", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D14", "D15", "D16", "D17", "D19", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D48", "D49", "D50", "D51")
foreach ($addr in $textFormatCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Updated Price (D) and Volume(1h) (E) values scraped for this run.
$ws.Range("D2").Value = '29.208.43'
$ws.Range("E2").Value = '  -0.04%  '
$ws.Range("D3").Value = '1.848.27'
$ws.Range("E3").Value = '  -0.65%  '
$ws.Range("D4").Value = '0.9993'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '246.39'
$ws.Range("D6").Value = '0.6982'
$ws.Range("E6").Value = '  -1.61%  '
$ws.Range("D7").Value = '0.9998'
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").Value = '0.07723'
$ws.Range("E8").Value = '  -1.10%  '
$ws.Range("D9").Value = '0.3059'
$ws.Range("E9").Value = '  -1.66%  '
$ws.Range("D10").Value = '23.53'
$ws.Range("E10").Value = '  -1.41%  '
$ws.Range("D11").Value = '0.07821'
$ws.Range("E11").Value = '  +0.05%  '
$ws.Range("D12").Value = '93.39'
$ws.Range("E12").Value = '  +0.80%  '
$ws.Range("D13").Value = '1.843.36'
$ws.Range("E13").Value = '  -0.10%  '
$ws.Range("D14").Value = '5.129'
$ws.Range("E14").Value = '  +0.04%  '
$ws.Range("D15").Value = '0.6860'
$ws.Range("E15").Value = '  -0.60%  '
$ws.Range("D16").Value = '6.635'
$ws.Range("E16").Value = '  +1.23%  '
$ws.Range("D17").Value = '0.000008313'
$ws.Range("E17").Value = '  -1.77%  '
$ws.Range("D18").Value = '29.199.13'
$ws.Range("E18").Value = '  -0.02%  '
$ws.Range("D19").Value = '241.28'
$ws.Range("E19").Value = '  -3.63%  '
$ws.Range("D20").Value = '2.083.90'
$ws.Range("E20").Value = '  -0.69%  '
$ws.Range("D21").Value = '12.76'
$ws.Range("E21").Value = '  -1.15%  '
$ws.Range("D22").Value = '0.9996'
$ws.Range("E22").Value = '  -0.01%  '
$ws.Range("D23").Value = '7.519'
$ws.Range("E23").Value = '  -0.42%  '
$ws.Range("D24").Value = '1.001'
$ws.Range("E24").Value = '  +0.07%  '
$ws.Range("D25").Value = '0.1516'
$ws.Range("E25").Value = '  -1.66%  '
$ws.Range("D26").Value = '159.09'
$ws.Range("E26").Value = '  -0.66%  '
$ws.Range("D27").Value = '8.825'
$ws.Range("E27").Value = '  -0.76%  '
$ws.Range("E28").Value = '  -1.64%  '
$ws.Range("D29").Value = '1.545'
$ws.Range("E29").Value = '  -1.19%  '
$ws.Range("D30").Value = '4.235'
$ws.Range("E30").Value = '  -1.06%  '
$ws.Range("D31").Value = '4.188'
$ws.Range("E31").Value = '  -1.44%  '
$ws.Range("D32").Value = '1.195'
$ws.Range("E32").Value = '  -0.78%  '
$ws.Range("D33").Value = '0.05122'
$ws.Range("E33").Value = '  -1.72%  '
$ws.Range("D34").Value = '0.7931'
$ws.Range("E34").Value = '  +4.47%  '
$ws.Range("D35").Value = '1.872'
$ws.Range("E35").Value = '  +1.20%  '
$ws.Range("D36").Value = '1.149'
$ws.Range("E36").Value = '  -2.31%  '
$ws.Range("D37").Value = '2.693'
$ws.Range("E37").Value = '  -0.68%  '
$ws.Range("D38").Value = '1.313.26'
$ws.Range("E38").Value = '  +7.06%  '
$ws.Range("D39").Value = '0.01871'
$ws.Range("E39").Value = '  +0.44%  '
$ws.Range("D40").Value = '2.712'
$ws.Range("E40").Value = '  -0.70%  '
$ws.Range("D41").Value = '0.9487'
$ws.Range("E41").Value = '  +5.69%  '
$ws.Range("D42").Value = '6.061'
$ws.Range("E42").Value = '  +6.77%  '
$ws.Range("D43").Value = '107.55'
$ws.Range("E43").Value = '  -1.85%  '
$ws.Range("D44").Value = '0.9996'
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("D45").Value = '9.731'
$ws.Range("E45").Value = '  +1.95%  '
$ws.Range("E46").Value = '  -0.66%  '
$ws.Range("D47").Value = '1.985.89'
$ws.Range("E47").Value = '  -0.57%  '
$ws.Range("D48").Value = '0.5178'
$ws.Range("E48").Value = '  -0.15%  '
$ws.Range("D49").Value = '64.21'
$ws.Range("E49").Value = '  -1.70%  '
$ws.Range("D50").Value = '1.766'
$ws.Range("E50").Value = '  +0.40%  '
$ws.Range("D51").Value = '6.998'
$ws.Range("E51").Value = '  -0.63%  '
